$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("examples")
$ws.Range("A1").Value = "test"
